$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 is a new match row being filled in (previously blank placeholder cells).
# Column A holds a date formatted as text (dd/mm/yyyy), e.g. "07/08/2025".
# Force it to be treated as plain text rather than being auto-parsed into a
# date serial number, then reset the style back to Normal so no extra
# number-format style sticks to the cell.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "07/08/2025"
$ws.Range("A21").Style = "Normal"

$ws.Range("B21").Value = "San Lorenzo"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "Velez Sarsfield"
$ws.Range("F21").Value = "L"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0.24
$ws.Range("L21").Value = 1.79
$ws.Range("M21").Value = 6
$ws.Range("N21").Value = 22
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 6
